$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had data on row 1 (B1) and a blank, custom-formatted
# row 3 (A3/C3 with an underline-font row style). Push that old row 3 down
# to row 6 (inserting 3 blank rows above it) so its row-level formatting
# (s="1" customFormat="1") ends up on row 6, matching the new layout.
$ws.Rows("3:5").Insert()

# Row 1: drop the old shared string in B1 (leave the cell blank but keep it
# underline-styled), add "indigena" in A1 and "Nome" in C1.
$ws.Range("B1").ClearContents()
$ws.Range("B1").Font.Underline = $true

$ws.Range("C4").Value = "CEP"
$ws.Range("B4").Value = "96784-346"
$ws.Range("A4").Value = "Rua Brinco de Princesa"

$ws.Range("A5").Value = "a33"

$ws.Range("C6").ClearContents()
$ws.Range("A6").Value = "a##"
$ws.Range("A6").Font.Underline = $true

$ws.Range("A2").Value = "sp"
$ws.Range("A2").Font.Underline = $true

$ws.Range("C1").Value = "Nome"
$ws.Range("C1").Font.ThemeColor = 1
$ws.Range("C2").Value = "Nome"

$ws.Range("A1").Value = "indigena"

$ws.Range("B2").Value = 3840033
$ws.Range("B3").Value = 38400321
$ws.Range("B3").Font.Underline = $true
$ws.Range("B5").Value = 38400322

$ws.Range("D1:D10").Select()

Write-Output "Done"
